$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New engagement rows 28-33 (services refactor: IWizytaService / WizytaService split
# into controller-facing interfaces + empty unit test scaffolding).
# Copy formatting (date style, s="3") from the last populated row (27) down
# through the new rows so the H column keeps its date number format.
$ws.Range("H27").Copy()
$ws.Range("H28:H33").PasteSpecial(-4122)

$ws.Range("H28").Value = 45788
$ws.Range("I28").Value = "WizytaService.cs"
$ws.Range("J28").Value = 20

$ws.Range("H29").Value = 45788
$ws.Range("I29").Value = "WykonaneBadaniaService.cs"
$ws.Range("J29").Value = 35

$ws.Range("H30").Value = 45788
$ws.Range("I30").Value = "OsobaService.cs"
$ws.Range("J30").Value = 20

$ws.Range("H31").Value = 45788
$ws.Range("I31").Value = "IOsobaService.cs"
$ws.Range("J31").Value = 9

$ws.Range("H32").Value = 45788
$ws.Range("I32").Value = "IWizytaService"
$ws.Range("J32").Value = 3

$ws.Range("H33").Value = 45788
$ws.Range("I33").Value = "IWykonaneBadaniaService.cs"
$ws.Range("J33").Value = 6
